$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old K:L columns (Level, Course) that are being removed
$ws.Range("K1:L10").Clear()

# Row 1: collapse header row into a single merged title cell
$ws.Range("B1:J1").Clear()
$ws.Range("A1").Value = "Herald College Kathmandu"

# Row 2
$ws.Range("A2").Value = "SUN"
$ws.Range("B2").Value = "9:30-11:30"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = "5CS024"
$ws.Range("E2").Value = "Collaborative Development"
$ws.Range("F2").Value = "Tutorial"
$ws.Range("G2").Value = "Mr. Udaya Kandel"
$ws.Range("H2").Value = "L5CG6"
$ws.Range("I2").Value = "WLV"
$ws.Range("J2").Value = "SR-02 Bilston"

# Row 3
$ws.Range("A3").Value = "SUN"
$ws.Range("B3").Value = "12:30-15:00"
$ws.Range("C3").Value = 2.5
$ws.Range("D3").Value = "5CS022"
$ws.Range("E3").Value = "Human Computer Interaction"
$ws.Range("F3").Value = "Workshop"
$ws.Range("G3").Value = "Mr. Dipesh Shrestha"
$ws.Range("H3").Value = "L5CG6"
$ws.Range("I3").Value = "WLV"
$ws.Range("J3").Value = "SR-04 Crompton"

# Row 4
$ws.Range("A4").Value = "MON"
$ws.Range("B4").Value = "9:30-12:00"
$ws.Range("C4").Value = 2.5
$ws.Range("D4").Value = "5CS024"
$ws.Range("E4").Value = "Collaborative Development"
$ws.Range("F4").Value = "Workshop"
$ws.Range("G4").Value = "Mr. Udaya Kandel"
$ws.Range("H4").Value = "L5CG6"
$ws.Range("I4").Value = "WLV"
$ws.Range("J4").Value = "SR-04 Crompton"

# Row 5
$ws.Range("A5").Value = "TUE"
$ws.Range("B5").Value = "7:00-9:00"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "5CS022"
$ws.Range("E5").Value = "Human Computer Interaction"
$ws.Range("F5").Value = "Lecture"
$ws.Range("G5").Value = "Mr. Apurba Neupane"
$ws.Range("H5").Value = "L5CG(5+6+7+8)"
$ws.Range("I5").Value = "WLV"
$ws.Range("J5").Value = "LT-02 Telford"

# Row 6
$ws.Range("A6").Value = "TUE"
$ws.Range("B6").Value = "9:30-11:30"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = "5CS020"
$ws.Range("E6").Value = "Distributed and Cloud Systems Programming"
$ws.Range("F6").Value = "Lecture"
$ws.Range("G6").Value = "Mr. Sumanta Silwal"
$ws.Range("H6").Value = "L5CG(5+6+7+8)"
$ws.Range("I6").Value = "WLV"
$ws.Range("J6").Value = "LT-01 Wulfruna"

# Row 7
$ws.Range("A7").Value = "WED"
$ws.Range("B7").Value = "7:00-9:00"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = "5CS024"
$ws.Range("E7").Value = "Collaborative Development"
$ws.Range("F7").Value = "Lecture"
$ws.Range("G7").Value = "Mr. Raj Shrestha"
$ws.Range("H7").Value = "L5CG(5+6+7+8)"
$ws.Range("I7").Value = "WLV"
$ws.Range("J7").Value = "LT-02 Telford"

# Row 8
$ws.Range("A8").Value = "WED"
$ws.Range("B8").Value = "9:30-11:30"
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = "5CS022"
$ws.Range("E8").Value = "Human Computer Interaction"
$ws.Range("F8").Value = "Tutorial"
$ws.Range("G8").Value = "Mr. Dipesh Shrestha"
$ws.Range("H8").Value = "L5CG6"
$ws.Range("I8").Value = "WLV"
$ws.Range("J8").Value = "SR-02 Bilston"

# Row 9
$ws.Range("A9").Value = "THU"
$ws.Range("B9").Value = "13:00-15:00"
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = "5CS020"
$ws.Range("E9").Value = "Distributed and Cloud Systems Programming"
$ws.Range("F9").Value = "Tutorial"
$ws.Range("G9").Value = "Mr. Prabin Sapkota"
$ws.Range("H9").Value = "L5CG6"
$ws.Range("I9").Value = "WLV"
$ws.Range("J9").Value = "SR-02 Bilston"

# Row 10
$ws.Range("A10").Value = "FRI"
$ws.Range("B10").Value = "12:30-14:30"
$ws.Range("C10").Value = 2.5
$ws.Range("D10").Value = "5CS020"
$ws.Range("E10").Value = "Distributed and Cloud Systems Programming"
$ws.Range("F10").Value = "Workshop"
$ws.Range("G10").Value = "Mr. Prabin Sapkota"
$ws.Range("H10").Value = "L5CG6"
$ws.Range("I10").Value = "WLV"
$ws.Range("J10").Value = "Lab-01 Mander"
